$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 28 (last quarter row) with refreshed metrics
$ws.Range("C28").Value = 558
$ws.Range("D28").Value = 55
$ws.Range("E28").Value = 503
$ws.Range("F28").Value = 8.566978193146417
